$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '52.222.02'
$r.ClearFormats()
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  -0.01%  '
$r.ClearFormats()

$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '2.840.73'
$r.ClearFormats()
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  +1.81%  '
$r.ClearFormats()

$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  -0.03%  '
$r.ClearFormats()

$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '360.92'
$r.ClearFormats()
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  +5.96%  '
$r.ClearFormats()

$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '113.38'
$r.ClearFormats()
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  -2.74%  '
$r.ClearFormats()

$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.576'
$r.ClearFormats()
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  +4.37%  '
$r.ClearFormats()

$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.ClearFormats()
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -0.02%  '
$r.ClearFormats()

$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.610'
$r.ClearFormats()
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  +5.22%  '
$r.ClearFormats()

$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '41.61'
$r.ClearFormats()
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -0.91%  '
$r.ClearFormats()

$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  -0.10%  '
$r.ClearFormats()

$r = $ws.Range('B12')
$r.NumberFormat = "@"
$r.Value = 'TRON'
$r.ClearFormats()
$r = $ws.Range('C12')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$r.ClearFormats()
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '0.132'
$r.ClearFormats()
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  +1.09%  '
$r.ClearFormats()

$r = $ws.Range('B13')
$r.NumberFormat = "@"
$r.Value = 'Chainlink'
$r.ClearFormats()
$r = $ws.Range('C13')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r.ClearFormats()
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '20.04'
$r.ClearFormats()
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  -0.28%  '
$r.ClearFormats()

$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  +2.58%  '
$r.ClearFormats()

$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '3.290.53'
$r.ClearFormats()
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  +1.71%  '
$r.ClearFormats()

$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '2.853.63'
$r.ClearFormats()
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  +1.58%  '
$r.ClearFormats()

$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '0.909'
$r.ClearFormats()
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  +2.70%  '
$r.ClearFormats()

$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '52.167.88'
$r.ClearFormats()
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  +0.09%  '
$r.ClearFormats()

$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '7.61'
$r.ClearFormats()
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  +9.44%  '
$r.ClearFormats()

$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '3.15'
$r.ClearFormats()
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  -1.76%  '
$r.ClearFormats()

$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '13.55'
$r.ClearFormats()
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  +1.72%  '
$r.ClearFormats()

$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  +1.32%  '
$r.ClearFormats()

$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '70.43'
$r.ClearFormats()
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  +0.23%  '
$r.ClearFormats()

$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '268.38'
$r.ClearFormats()
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  -3.65%  '
$r.ClearFormats()

$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  +1.48%  '
$r.ClearFormats()

$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '27.14'
$r.ClearFormats()
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  +1.11%  '
$r.ClearFormats()

$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  +0.10%  '
$r.ClearFormats()

$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  +1.96%  '
$r.ClearFormats()

$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  +1.46%  '
$r.ClearFormats()

$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '54.04'
$r.ClearFormats()
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  +7.24%  '
$r.ClearFormats()

$r = $ws.Range('B31')
$r.NumberFormat = "@"
$r.Value = 'VeChain'
$r.ClearFormats()
$r = $ws.Range('C31')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r.ClearFormats()
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '0.0485'
$r.ClearFormats()
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  +29.29%  '
$r.ClearFormats()

$r = $ws.Range('B32')
$r.NumberFormat = "@"
$r.Value = 'Kaspa'
$r.ClearFormats()
$r = $ws.Range('C32')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$r.ClearFormats()
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '0.141'
$r.ClearFormats()
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  -1.05%  '
$r.ClearFormats()

$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '34.74'
$r.ClearFormats()
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -0.06%  '
$r.ClearFormats()

$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '5.89'
$r.ClearFormats()
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  +2.45%  '
$r.ClearFormats()

$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '5.44'
$r.ClearFormats()
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  +9.13%  '
$r.ClearFormats()

$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '0.0846'
$r.ClearFormats()
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  +2.23%  '
$r.ClearFormats()

$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  +0.01%  '
$r.ClearFormats()

$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '3.28'
$r.ClearFormats()
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  +1.10%  '
$r.ClearFormats()

$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  -2.18%  '
$r.ClearFormats()

$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '18.42'
$r.ClearFormats()
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  -2.80%  '
$r.ClearFormats()

$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '23.96'
$r.ClearFormats()
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  +2.06%  '
$r.ClearFormats()

$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  +1.48%  '
$r.ClearFormats()

$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '128.02'
$r.ClearFormats()
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  +2.67%  '
$r.ClearFormats()

$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -6.91%  '
$r.ClearFormats()

$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -1.99%  '
$r.ClearFormats()

$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  +3.15%  '
$r.ClearFormats()

$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '2.115.17'
$r.ClearFormats()
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  +0.82%  '
$r.ClearFormats()

$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  +1.06%  '
$r.ClearFormats()

$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  +11.03%  '
$r.ClearFormats()

$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '5.88'
$r.ClearFormats()
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  +5.44%  '
$r.ClearFormats()

$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '61.85'
$r.ClearFormats()
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  +2.57%  '
$r.ClearFormats()
